$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("FECHA","CANT","DESCRIPCION","TALLA","NOMBRE","TELEFONO","TIPO DE PAGO","A CUENTA","EFECTIVO","TARJETA","TOTAL DE VENTA","CIERRE/DIA","TOTAL")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(2, $col).Value = $headers[$i]
}

$ws.Range("A1:A2").Borders.LineStyle = 1

$ws.Range("M3").Select()
